# Carter Kreis timesheet: add the 10/24 work-session entries (rows 10-11)
# and point the saved selection at the newly-entered data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 -----------------------------------------------------------
$ws.Range("A10").Value = 45954                         # 10/24/2025
$ws.Range("B10").Value = "Pre-Alpha Build"
$ws.Range("C10").Value = 0.10416666666666667            # 2:30 AM -> h:mm
$ws.Range("C10").NumberFormat = "h:mm"
$ws.Range("D10").Value = 0.1701388888888889              # 4:05 AM -> h:mm
$ws.Range("D10").NumberFormat = "h:mm"
$ws.Range("F10").Value = "Setup ESP-IDF extension in VS Code. Attempted Hello World example code but my laptop is not detecting ESP32 as a COM Port. We also continued to plan Pre-Alpha Build goal"

# --- Row 11 -----------------------------------------------------------
$ws.Range("A11").Value = 45954                         # 10/24/2025
$ws.Range("A11").NumberFormat = "m/d/yy"                 # same date style as A10
$ws.Range("B11").Value = "Pre-Alpha Build"
$ws.Range("C11").Value = 0.30208333333333331
$ws.Range("C11").NumberFormat = "h:mm"
$ws.Range("D11").Value = 0.37361111111111112
$ws.Range("D11").NumberFormat = "h:mm"
$ws.Range("F11").Value = "Repeated issues with recursively cloning the submodule esp-idf. Cloning into the submodule itself took 10-12 minutes and the first two times some of the directories failed to clone on time."

# --- View / selection ---------------------------------------------------
$ws.Activate()
$ws.Range("C12").Select()
